$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 3).Value = 6
$ws.Cells.Item(2, 6).Value = 82
$ws.Cells.Item(2, 12).Value = 'stimuli/img_o30wb.png'
$ws.Cells.Item(2, 13).Value = 81.06666666666666
$ws.Cells.Item(2, 14).Value = 65.37777777777778
$ws.Cells.Item(2, 15).Value = 73.22222222222223
$ws.Cells.Item(2, 16).Value = 45
$ws.Cells.Item(2, 17).Value = 8
$ws.Cells.Item(2, 18).Value = 8
$ws.Cells.Item(2, 19).Value = 8
$ws.Cells.Item(2, 20).Value = 8
$ws.Cells.Item(2, 21).Value = 8
$ws.Cells.Item(2, 22).Value = 8

# Row 3
$ws.Cells.Item(3, 3).Value = 6
$ws.Cells.Item(3, 6).Value = 83
$ws.Cells.Item(3, 8).Value = 'bedrooms'
$ws.Cells.Item(3, 9).Value = 'distractor'
$ws.Cells.Item(3, 11).Value = 'f'
$ws.Cells.Item(3, 12).Value = 'stimuli/img_7caxh.png'
$ws.Cells.Item(3, 13).Value = 83.38095238095238
$ws.Cells.Item(3, 14).Value = 64.26190476190476
$ws.Cells.Item(3, 15).Value = 73.82142857142857
$ws.Cells.Item(3, 17).Value = 8
$ws.Cells.Item(3, 18).Value = 8
$ws.Cells.Item(3, 19).Value = 8
$ws.Cells.Item(3, 20).Value = 8
$ws.Cells.Item(3, 21).Value = 8
$ws.Cells.Item(3, 22).Value = 8

# Row 4
$ws.Cells.Item(4, 3).Value = 6
$ws.Cells.Item(4, 6).Value = 84
$ws.Cells.Item(4, 8).Value = 'kitchens'
$ws.Cells.Item(4, 9).Value = 'distractor'
$ws.Cells.Item(4, 11).Value = 'f'
$ws.Cells.Item(4, 12).Value = 'stimuli/img_lgyo6.png'
$ws.Cells.Item(4, 13).Value = 54.68421052631579
$ws.Cells.Item(4, 14).Value = 27.18421052631579
$ws.Cells.Item(4, 15).Value = 40.93421052631579
$ws.Cells.Item(4, 16).Value = 38
$ws.Cells.Item(4, 17).Value = 2

# Row 5
$ws.Cells.Item(5, 3).Value = 6
$ws.Cells.Item(5, 6).Value = 85
$ws.Cells.Item(5, 8).Value = 'bedrooms'
$ws.Cells.Item(5, 9).Value = 'distractor'
$ws.Cells.Item(5, 11).Value = 'f'
$ws.Cells.Item(5, 12).Value = 'stimuli/img_91csq.png'
$ws.Cells.Item(5, 13).Value = 50.44736842105263
$ws.Cells.Item(5, 14).Value = 28.34210526315789
$ws.Cells.Item(5, 15).Value = 39.39473684210526
$ws.Cells.Item(5, 16).Value = 38
$ws.Cells.Item(5, 17).Value = 2
$ws.Cells.Item(5, 18).Value = 2
$ws.Cells.Item(5, 19).Value = 2
$ws.Cells.Item(5, 20).Value = 2
$ws.Cells.Item(5, 21).Value = 2
$ws.Cells.Item(5, 22).Value = 2

# Row 6
$ws.Cells.Item(6, 3).Value = 6
$ws.Cells.Item(6, 6).Value = 86
$ws.Cells.Item(6, 12).Value = 'stimuli/img_a9he3.png'
$ws.Cells.Item(6, 13).Value = 83.06521739130434
$ws.Cells.Item(6, 14).Value = 63.95652173913044
$ws.Cells.Item(6, 15).Value = 73.51086956521739
$ws.Cells.Item(6, 17).Value = 8
$ws.Cells.Item(6, 18).Value = 8
$ws.Cells.Item(6, 19).Value = 8
$ws.Cells.Item(6, 20).Value = 8
$ws.Cells.Item(6, 21).Value = 8
$ws.Cells.Item(6, 22).Value = 8

# Row 7
$ws.Cells.Item(7, 3).Value = 6
$ws.Cells.Item(7, 6).Value = 87
$ws.Cells.Item(7, 8).Value = 'living_rooms'
$ws.Cells.Item(7, 9).Value = 'target'
$ws.Cells.Item(7, 11).Value = 'j'
$ws.Cells.Item(7, 12).Value = 'stimuli/img_9bkl9.png'
$ws.Cells.Item(7, 13).Value = 46.62162162162162
$ws.Cells.Item(7, 14).Value = 34.27027027027027
$ws.Cells.Item(7, 15).Value = 40.44594594594595
$ws.Cells.Item(7, 16).Value = 37

# Row 8
$ws.Cells.Item(8, 3).Value = 6
$ws.Cells.Item(8, 6).Value = 88
$ws.Cells.Item(8, 8).Value = 'kitchens'
$ws.Cells.Item(8, 12).Value = 'stimuli/img_xtb79.png'
$ws.Cells.Item(8, 13).Value = 55.59375
$ws.Cells.Item(8, 14).Value = 31.40625
$ws.Cells.Item(8, 15).Value = 43.5
$ws.Cells.Item(8, 16).Value = 32
$ws.Cells.Item(8, 17).Value = 2
$ws.Cells.Item(8, 18).Value = 2
$ws.Cells.Item(8, 19).Value = 2
$ws.Cells.Item(8, 20).Value = 2
$ws.Cells.Item(8, 21).Value = 2
$ws.Cells.Item(8, 22).Value = 2

# Row 9
$ws.Cells.Item(9, 3).Value = 6
$ws.Cells.Item(9, 6).Value = 89
$ws.Cells.Item(9, 12).Value = 'stimuli/img_wgddx.png'
$ws.Cells.Item(9, 13).Value = 45.6304347826087
$ws.Cells.Item(9, 14).Value = 34.30434782608695
$ws.Cells.Item(9, 15).Value = 39.96739130434783
$ws.Cells.Item(9, 16).Value = 46
$ws.Cells.Item(9, 17).Value = 3
$ws.Cells.Item(9, 18).Value = 3
$ws.Cells.Item(9, 19).Value = 3
$ws.Cells.Item(9, 20).Value = 3
$ws.Cells.Item(9, 21).Value = 3
$ws.Cells.Item(9, 22).Value = 4

# Row 10
$ws.Cells.Item(10, 3).Value = 6
$ws.Cells.Item(10, 6).Value = 90
$ws.Cells.Item(10, 12).Value = 'stimuli/img_zxvl3.png'
$ws.Cells.Item(10, 13).Value = 68.78260869565217
$ws.Cells.Item(10, 14).Value = 47.56521739130435
$ws.Cells.Item(10, 15).Value = 58.17391304347827
$ws.Cells.Item(10, 16).Value = 46
$ws.Cells.Item(10, 17).Value = 5
$ws.Cells.Item(10, 18).Value = 5
$ws.Cells.Item(10, 19).Value = 5
$ws.Cells.Item(10, 20).Value = 5
$ws.Cells.Item(10, 21).Value = 5
$ws.Cells.Item(10, 22).Value = 5

# Row 11
$ws.Cells.Item(11, 3).Value = 6
$ws.Cells.Item(11, 6).Value = 91
$ws.Cells.Item(11, 8).Value = 'bedrooms'
$ws.Cells.Item(11, 12).Value = 'stimuli/img_th7xh.png'
$ws.Cells.Item(11, 13).Value = 82.35897435897436
$ws.Cells.Item(11, 14).Value = 65.53846153846153
$ws.Cells.Item(11, 15).Value = 73.94871794871796
$ws.Cells.Item(11, 16).Value = 39
$ws.Cells.Item(11, 17).Value = 8
$ws.Cells.Item(11, 18).Value = 8
$ws.Cells.Item(11, 19).Value = 8
$ws.Cells.Item(11, 20).Value = 8
$ws.Cells.Item(11, 21).Value = 8
$ws.Cells.Item(11, 22).Value = 8

# Row 12
$ws.Cells.Item(12, 3).Value = 6
$ws.Cells.Item(12, 6).Value = 92
$ws.Cells.Item(12, 8).Value = 'bedrooms'
$ws.Cells.Item(12, 9).Value = 'distractor'
$ws.Cells.Item(12, 11).Value = 'f'
$ws.Cells.Item(12, 12).Value = 'stimuli/img_ys3qz.png'
$ws.Cells.Item(12, 13).Value = 46.79545454545455
$ws.Cells.Item(12, 14).Value = 31.20454545454545
$ws.Cells.Item(12, 15).Value = 39
$ws.Cells.Item(12, 16).Value = 44
$ws.Cells.Item(12, 17).Value = 2
$ws.Cells.Item(12, 18).Value = 2
$ws.Cells.Item(12, 19).Value = 2
$ws.Cells.Item(12, 20).Value = 2
$ws.Cells.Item(12, 21).Value = 2
$ws.Cells.Item(12, 22).Value = 2

# Row 13
$ws.Cells.Item(13, 3).Value = 6
$ws.Cells.Item(13, 6).Value = 93
$ws.Cells.Item(13, 8).Value = 'living_rooms'
$ws.Cells.Item(13, 9).Value = 'target'
$ws.Cells.Item(13, 11).Value = 'j'
$ws.Cells.Item(13, 12).Value = 'stimuli/img_c0vzo.png'
$ws.Cells.Item(13, 13).Value = 21.51162790697675
$ws.Cells.Item(13, 14).Value = 8.232558139534884
$ws.Cells.Item(13, 15).Value = 14.87209302325581
$ws.Cells.Item(13, 16).Value = 43
$ws.Cells.Item(13, 17).Value = 1
$ws.Cells.Item(13, 18).Value = 1
$ws.Cells.Item(13, 19).Value = 1
$ws.Cells.Item(13, 20).Value = 1
$ws.Cells.Item(13, 21).Value = 1
$ws.Cells.Item(13, 22).Value = 1

# Row 14
$ws.Cells.Item(14, 3).Value = 6
$ws.Cells.Item(14, 6).Value = 94
$ws.Cells.Item(14, 12).Value = 'stimuli/img_rych7.png'
$ws.Cells.Item(14, 13).Value = 30.4468085106383
$ws.Cells.Item(14, 14).Value = 23.4468085106383
$ws.Cells.Item(14, 15).Value = 26.9468085106383
$ws.Cells.Item(14, 16).Value = 47
$ws.Cells.Item(14, 17).Value = 2
$ws.Cells.Item(14, 18).Value = 2
$ws.Cells.Item(14, 19).Value = 2
$ws.Cells.Item(14, 20).Value = 2
$ws.Cells.Item(14, 21).Value = 2
$ws.Cells.Item(14, 22).Value = 2

# Row 15
$ws.Cells.Item(15, 3).Value = 6
$ws.Cells.Item(15, 6).Value = 95
$ws.Cells.Item(15, 8).Value = 'kitchens'
$ws.Cells.Item(15, 9).Value = 'distractor'
$ws.Cells.Item(15, 11).Value = 'f'
$ws.Cells.Item(15, 12).Value = 'stimuli/img_uhmld.png'
$ws.Cells.Item(15, 13).Value = 54.19354838709678
$ws.Cells.Item(15, 14).Value = 32
$ws.Cells.Item(15, 15).Value = 43.09677419354838
$ws.Cells.Item(15, 16).Value = 31
$ws.Cells.Item(15, 17).Value = 2
$ws.Cells.Item(15, 18).Value = 2
$ws.Cells.Item(15, 19).Value = 2
$ws.Cells.Item(15, 20).Value = 2
$ws.Cells.Item(15, 21).Value = 2
$ws.Cells.Item(15, 22).Value = 2

# Row 16
$ws.Cells.Item(16, 3).Value = 6
$ws.Cells.Item(16, 6).Value = 96
$ws.Cells.Item(16, 8).Value = 'living_rooms'
$ws.Cells.Item(16, 9).Value = 'target'
$ws.Cells.Item(16, 11).Value = 'j'
$ws.Cells.Item(16, 12).Value = 'stimuli/img_dg5h7.png'
$ws.Cells.Item(16, 13).Value = 88.72093023255815
$ws.Cells.Item(16, 14).Value = 76.06976744186046
$ws.Cells.Item(16, 15).Value = 82.3953488372093
$ws.Cells.Item(16, 16).Value = 43
$ws.Cells.Item(16, 17).Value = 10
$ws.Cells.Item(16, 18).Value = 10
$ws.Cells.Item(16, 19).Value = 10
$ws.Cells.Item(16, 20).Value = 10
$ws.Cells.Item(16, 21).Value = 10
$ws.Cells.Item(16, 22).Value = 10

# Row 17
$ws.Cells.Item(17, 3).Value = 6
$ws.Cells.Item(17, 6).Value = 97
$ws.Cells.Item(17, 8).Value = 'living_rooms'
$ws.Cells.Item(17, 9).Value = 'target'
$ws.Cells.Item(17, 11).Value = 'j'
$ws.Cells.Item(17, 12).Value = 'stimuli/img_i6wsx.png'
$ws.Cells.Item(17, 13).Value = 79.07142857142857
$ws.Cells.Item(17, 14).Value = 58
$ws.Cells.Item(17, 15).Value = 68.53571428571428
$ws.Cells.Item(17, 16).Value = 42
$ws.Cells.Item(17, 17).Value = 7
$ws.Cells.Item(17, 18).Value = 7
$ws.Cells.Item(17, 19).Value = 7
$ws.Cells.Item(17, 20).Value = 7
$ws.Cells.Item(17, 21).Value = 7
$ws.Cells.Item(17, 22).Value = 7

# Row 18
$ws.Cells.Item(18, 3).Value = 6
$ws.Cells.Item(18, 6).Value = 98
$ws.Cells.Item(18, 8).Value = 'living_rooms'
$ws.Cells.Item(18, 9).Value = 'target'
$ws.Cells.Item(18, 11).Value = 'j'
$ws.Cells.Item(18, 12).Value = 'stimuli/img_jkm86.png'
$ws.Cells.Item(18, 13).Value = 58.32558139534883
$ws.Cells.Item(18, 14).Value = 38.65116279069768
$ws.Cells.Item(18, 15).Value = 48.48837209302326
$ws.Cells.Item(18, 16).Value = 43
$ws.Cells.Item(18, 17).Value = 4
$ws.Cells.Item(18, 18).Value = 4
$ws.Cells.Item(18, 19).Value = 4
$ws.Cells.Item(18, 20).Value = 4
$ws.Cells.Item(18, 21).Value = 4
$ws.Cells.Item(18, 22).Value = 4

# Row 19
$ws.Cells.Item(19, 3).Value = 6
$ws.Cells.Item(19, 6).Value = 99
$ws.Cells.Item(19, 8).Value = 'kitchens'
$ws.Cells.Item(19, 12).Value = 'stimuli/img_c30d1.png'
$ws.Cells.Item(19, 13).Value = 78.875
$ws.Cells.Item(19, 14).Value = 60.34375
$ws.Cells.Item(19, 15).Value = 69.609375
$ws.Cells.Item(19, 16).Value = 32
$ws.Cells.Item(19, 17).Value = 8
$ws.Cells.Item(19, 18).Value = 8
$ws.Cells.Item(19, 19).Value = 8
$ws.Cells.Item(19, 20).Value = 8
$ws.Cells.Item(19, 21).Value = 8
$ws.Cells.Item(19, 22).Value = 8

# Row 20
$ws.Cells.Item(20, 3).Value = 6
$ws.Cells.Item(20, 6).Value = 100
$ws.Cells.Item(20, 12).Value = 'stimuli/img_xzyzy.png'
$ws.Cells.Item(20, 13).Value = 85.37209302325581
$ws.Cells.Item(20, 14).Value = 68.90697674418605
$ws.Cells.Item(20, 15).Value = 77.13953488372093
$ws.Cells.Item(20, 16).Value = 43
$ws.Cells.Item(20, 17).Value = 9
$ws.Cells.Item(20, 18).Value = 9
$ws.Cells.Item(20, 19).Value = 9
$ws.Cells.Item(20, 20).Value = 9
$ws.Cells.Item(20, 21).Value = 9
$ws.Cells.Item(20, 22).Value = 9

# Row 21
$ws.Cells.Item(21, 3).Value = 6
$ws.Cells.Item(21, 6).Value = 101
$ws.Cells.Item(21, 8).Value = 'kitchens'
$ws.Cells.Item(21, 12).Value = 'stimuli/img_3tnh4.png'
$ws.Cells.Item(21, 13).Value = 80.43243243243244
$ws.Cells.Item(21, 14).Value = 58.72972972972973
$ws.Cells.Item(21, 15).Value = 69.58108108108108
$ws.Cells.Item(21, 16).Value = 37
$ws.Cells.Item(21, 17).Value = 8
$ws.Cells.Item(21, 18).Value = 8
$ws.Cells.Item(21, 19).Value = 8
$ws.Cells.Item(21, 20).Value = 8
$ws.Cells.Item(21, 21).Value = 8
$ws.Cells.Item(21, 22).Value = 8

# Row 22
$ws.Cells.Item(22, 3).Value = 6
$ws.Cells.Item(22, 6).Value = 102
$ws.Cells.Item(22, 12).Value = 'stimuli/img_gztbt.png'
$ws.Cells.Item(22, 13).Value = 55.06451612903226
$ws.Cells.Item(22, 14).Value = 26.09677419354839
$ws.Cells.Item(22, 15).Value = 40.58064516129032
$ws.Cells.Item(22, 16).Value = 31
$ws.Cells.Item(22, 17).Value = 2
$ws.Cells.Item(22, 18).Value = 2
$ws.Cells.Item(22, 19).Value = 2
$ws.Cells.Item(22, 20).Value = 2
$ws.Cells.Item(22, 21).Value = 2
$ws.Cells.Item(22, 22).Value = 2

# Row 23
$ws.Cells.Item(23, 3).Value = 6
$ws.Cells.Item(23, 6).Value = 103
$ws.Cells.Item(23, 8).Value = 'bedrooms'
$ws.Cells.Item(23, 9).Value = 'distractor'
$ws.Cells.Item(23, 11).Value = 'f'
$ws.Cells.Item(23, 12).Value = 'stimuli/img_swcci.png'
$ws.Cells.Item(23, 13).Value = 49.82926829268293
$ws.Cells.Item(23, 14).Value = 28.46341463414634
$ws.Cells.Item(23, 15).Value = 39.14634146341464
$ws.Cells.Item(23, 16).Value = 41
$ws.Cells.Item(23, 17).Value = 2
$ws.Cells.Item(23, 18).Value = 2
$ws.Cells.Item(23, 19).Value = 2
$ws.Cells.Item(23, 20).Value = 2
$ws.Cells.Item(23, 21).Value = 2
$ws.Cells.Item(23, 22).Value = 2

# Row 24
$ws.Cells.Item(24, 3).Value = 6
$ws.Cells.Item(24, 6).Value = 104
$ws.Cells.Item(24, 12).Value = 'stimuli/img_g13d5.png'
$ws.Cells.Item(24, 13).Value = 73
$ws.Cells.Item(24, 14).Value = 51.51111111111111
$ws.Cells.Item(24, 15).Value = 62.25555555555556
$ws.Cells.Item(24, 16).Value = 45
$ws.Cells.Item(24, 17).Value = 6
$ws.Cells.Item(24, 18).Value = 6
$ws.Cells.Item(24, 19).Value = 6
$ws.Cells.Item(24, 20).Value = 6
$ws.Cells.Item(24, 21).Value = 6
$ws.Cells.Item(24, 22).Value = 6

# Row 25
$ws.Cells.Item(25, 3).Value = 6
$ws.Cells.Item(25, 6).Value = 105
$ws.Cells.Item(25, 12).Value = 'stimuli/img_7wquy.png'
$ws.Cells.Item(25, 13).Value = 50.59375
$ws.Cells.Item(25, 14).Value = 30.59375
$ws.Cells.Item(25, 15).Value = 40.59375
$ws.Cells.Item(25, 16).Value = 32
$ws.Cells.Item(25, 17).Value = 2
$ws.Cells.Item(25, 18).Value = 2
$ws.Cells.Item(25, 19).Value = 2
$ws.Cells.Item(25, 20).Value = 2
$ws.Cells.Item(25, 21).Value = 2
$ws.Cells.Item(25, 22).Value = 2

# Row 26
$ws.Cells.Item(26, 3).Value = 6
$ws.Cells.Item(26, 6).Value = 106
$ws.Cells.Item(26, 8).Value = 'living_rooms'
$ws.Cells.Item(26, 9).Value = 'target'
$ws.Cells.Item(26, 11).Value = 'j'
$ws.Cells.Item(26, 12).Value = 'stimuli/img_165pk.png'
$ws.Cells.Item(26, 13).Value = 85.73333333333333
$ws.Cells.Item(26, 14).Value = 69.22222222222223
$ws.Cells.Item(26, 15).Value = 77.47777777777779
$ws.Cells.Item(26, 16).Value = 45
$ws.Cells.Item(26, 17).Value = 9
$ws.Cells.Item(26, 18).Value = 9
$ws.Cells.Item(26, 19).Value = 9
$ws.Cells.Item(26, 20).Value = 9
$ws.Cells.Item(26, 21).Value = 9
$ws.Cells.Item(26, 22).Value = 9

# Row 27
$ws.Cells.Item(27, 3).Value = 6
$ws.Cells.Item(27, 6).Value = 107
$ws.Cells.Item(27, 12).Value = 'stimuli/img_ymgcb.png'
$ws.Cells.Item(27, 13).Value = 83.73684210526316
$ws.Cells.Item(27, 14).Value = 61.13157894736842
$ws.Cells.Item(27, 15).Value = 72.4342105263158
$ws.Cells.Item(27, 16).Value = 38
$ws.Cells.Item(27, 17).Value = 8
$ws.Cells.Item(27, 18).Value = 8
$ws.Cells.Item(27, 19).Value = 8
$ws.Cells.Item(27, 20).Value = 8
$ws.Cells.Item(27, 21).Value = 8
$ws.Cells.Item(27, 22).Value = 8

# Row 28
$ws.Cells.Item(28, 3).Value = 6
$ws.Cells.Item(28, 6).Value = 108
$ws.Cells.Item(28, 8).Value = 'living_rooms'
$ws.Cells.Item(28, 9).Value = 'target'
$ws.Cells.Item(28, 11).Value = 'j'
$ws.Cells.Item(28, 12).Value = 'stimuli/img_5jy9c.png'
$ws.Cells.Item(28, 13).Value = 87.37209302325581
$ws.Cells.Item(28, 14).Value = 79.18604651162791
$ws.Cells.Item(28, 15).Value = 83.27906976744185
$ws.Cells.Item(28, 16).Value = 43
$ws.Cells.Item(28, 17).Value = 10
$ws.Cells.Item(28, 18).Value = 10
$ws.Cells.Item(28, 19).Value = 10
$ws.Cells.Item(28, 20).Value = 10
$ws.Cells.Item(28, 21).Value = 9
$ws.Cells.Item(28, 22).Value = 10

# Row 29
$ws.Cells.Item(29, 3).Value = 6
$ws.Cells.Item(29, 6).Value = 109
$ws.Cells.Item(29, 8).Value = 'bedrooms'
$ws.Cells.Item(29, 9).Value = 'distractor'
$ws.Cells.Item(29, 11).Value = 'f'
$ws.Cells.Item(29, 12).Value = 'stimuli/img_gv750.png'
$ws.Cells.Item(29, 13).Value = 83.51428571428572
$ws.Cells.Item(29, 14).Value = 61.88571428571429
$ws.Cells.Item(29, 15).Value = 72.7
$ws.Cells.Item(29, 16).Value = 35
$ws.Cells.Item(29, 17).Value = 8
$ws.Cells.Item(29, 18).Value = 8
$ws.Cells.Item(29, 19).Value = 8
$ws.Cells.Item(29, 20).Value = 8
$ws.Cells.Item(29, 21).Value = 8
$ws.Cells.Item(29, 22).Value = 8

# Row 30
$ws.Cells.Item(30, 3).Value = 6
$ws.Cells.Item(30, 6).Value = 110
$ws.Cells.Item(30, 12).Value = 'stimuli/img_w8yhd.png'
$ws.Cells.Item(30, 13).Value = 55.74418604651163
$ws.Cells.Item(30, 14).Value = 38.90697674418605
$ws.Cells.Item(30, 15).Value = 47.32558139534883
$ws.Cells.Item(30, 16).Value = 43
$ws.Cells.Item(30, 17).Value = 4
$ws.Cells.Item(30, 18).Value = 4
$ws.Cells.Item(30, 19).Value = 4
$ws.Cells.Item(30, 20).Value = 4
$ws.Cells.Item(30, 21).Value = 4
$ws.Cells.Item(30, 22).Value = 4

# Row 31
$ws.Cells.Item(31, 3).Value = 6
$ws.Cells.Item(31, 6).Value = 111
$ws.Cells.Item(31, 8).Value = 'living_rooms'
$ws.Cells.Item(31, 9).Value = 'target'
$ws.Cells.Item(31, 11).Value = 'j'
$ws.Cells.Item(31, 12).Value = 'stimuli/img_nb8p4.png'
$ws.Cells.Item(31, 13).Value = 16.36170212765957
$ws.Cells.Item(31, 14).Value = 12.70212765957447
$ws.Cells.Item(31, 15).Value = 14.53191489361702
$ws.Cells.Item(31, 16).Value = 47
$ws.Cells.Item(31, 17).Value = 1
$ws.Cells.Item(31, 18).Value = 1
$ws.Cells.Item(31, 19).Value = 1
$ws.Cells.Item(31, 20).Value = 1
$ws.Cells.Item(31, 21).Value = 1
$ws.Cells.Item(31, 22).Value = 1

# Row 32
$ws.Cells.Item(32, 3).Value = 6
$ws.Cells.Item(32, 6).Value = 112
$ws.Cells.Item(32, 8).Value = 'kitchens'
$ws.Cells.Item(32, 9).Value = 'distractor'
$ws.Cells.Item(32, 11).Value = 'f'
$ws.Cells.Item(32, 12).Value = 'stimuli/img_463mq.png'
$ws.Cells.Item(32, 13).Value = 51.35294117647059
$ws.Cells.Item(32, 14).Value = 30.20588235294118
$ws.Cells.Item(32, 15).Value = 40.77941176470588
$ws.Cells.Item(32, 16).Value = 34
$ws.Cells.Item(32, 17).Value = 2
$ws.Cells.Item(32, 18).Value = 2
$ws.Cells.Item(32, 19).Value = 2
$ws.Cells.Item(32, 20).Value = 2
$ws.Cells.Item(32, 21).Value = 2
$ws.Cells.Item(32, 22).Value = 2

# Row 33
$ws.Cells.Item(33, 3).Value = 6
$ws.Cells.Item(33, 6).Value = 113
$ws.Cells.Item(33, 12).Value = 'stimuli/img_jpldg.png'
$ws.Cells.Item(33, 13).Value = 79.54545454545455
$ws.Cells.Item(33, 14).Value = 57.75
$ws.Cells.Item(33, 15).Value = 68.64772727272728
$ws.Cells.Item(33, 16).Value = 44
$ws.Cells.Item(33, 17).Value = 7
$ws.Cells.Item(33, 18).Value = 7
$ws.Cells.Item(33, 19).Value = 7
$ws.Cells.Item(33, 20).Value = 7
$ws.Cells.Item(33, 21).Value = 7
$ws.Cells.Item(33, 22).Value = 7

# Row 34
$ws.Cells.Item(34, 3).Value = 6
$ws.Cells.Item(34, 6).Value = 114
$ws.Cells.Item(34, 8).Value = 'bedrooms'
$ws.Cells.Item(34, 12).Value = 'stimuli/img_md6k3.png'
$ws.Cells.Item(34, 13).Value = 81.57142857142857
$ws.Cells.Item(34, 14).Value = 63.66666666666666
$ws.Cells.Item(34, 15).Value = 72.61904761904762
$ws.Cells.Item(34, 16).Value = 42
$ws.Cells.Item(34, 17).Value = 8
$ws.Cells.Item(34, 18).Value = 8
$ws.Cells.Item(34, 19).Value = 8
$ws.Cells.Item(34, 20).Value = 8
$ws.Cells.Item(34, 21).Value = 8
$ws.Cells.Item(34, 22).Value = 8

# Row 35
$ws.Cells.Item(35, 3).Value = 6
$ws.Cells.Item(35, 6).Value = 115
$ws.Cells.Item(35, 8).Value = 'bedrooms'
$ws.Cells.Item(35, 9).Value = 'distractor'
$ws.Cells.Item(35, 11).Value = 'f'
$ws.Cells.Item(35, 12).Value = 'stimuli/img_awuev.png'
$ws.Cells.Item(35, 13).Value = 44.21052631578947
$ws.Cells.Item(35, 14).Value = 24.26315789473684
$ws.Cells.Item(35, 15).Value = 34.23684210526316
$ws.Cells.Item(35, 16).Value = 38
$ws.Cells.Item(35, 17).Value = 2
$ws.Cells.Item(35, 18).Value = 2
$ws.Cells.Item(35, 19).Value = 2
$ws.Cells.Item(35, 20).Value = 2
$ws.Cells.Item(35, 21).Value = 2
$ws.Cells.Item(35, 22).Value = 2

# Row 36
$ws.Cells.Item(36, 3).Value = 6
$ws.Cells.Item(36, 6).Value = 116
$ws.Cells.Item(36, 8).Value = 'bedrooms'
$ws.Cells.Item(36, 9).Value = 'distractor'
$ws.Cells.Item(36, 11).Value = 'f'
$ws.Cells.Item(36, 12).Value = 'stimuli/img_b2jli.png'
$ws.Cells.Item(36, 13).Value = 83.15625
$ws.Cells.Item(36, 14).Value = 63.8125
$ws.Cells.Item(36, 15).Value = 73.484375
$ws.Cells.Item(36, 16).Value = 32
$ws.Cells.Item(36, 17).Value = 8
$ws.Cells.Item(36, 18).Value = 8
$ws.Cells.Item(36, 19).Value = 8
$ws.Cells.Item(36, 20).Value = 8
$ws.Cells.Item(36, 21).Value = 8
$ws.Cells.Item(36, 22).Value = 8

# Row 37
$ws.Cells.Item(37, 3).Value = 6
$ws.Cells.Item(37, 6).Value = 117
$ws.Cells.Item(37, 8).Value = 'living_rooms'
$ws.Cells.Item(37, 9).Value = 'target'
$ws.Cells.Item(37, 11).Value = 'j'
$ws.Cells.Item(37, 12).Value = 'stimuli/img_8dmpq.png'
$ws.Cells.Item(37, 13).Value = 30.65909090909091
$ws.Cells.Item(37, 14).Value = 24.11363636363636
$ws.Cells.Item(37, 15).Value = 27.38636363636364
$ws.Cells.Item(37, 16).Value = 44
$ws.Cells.Item(37, 17).Value = 2
$ws.Cells.Item(37, 18).Value = 2
$ws.Cells.Item(37, 19).Value = 2
$ws.Cells.Item(37, 20).Value = 2
$ws.Cells.Item(37, 21).Value = 2
$ws.Cells.Item(37, 22).Value = 2

# Row 38
$ws.Cells.Item(38, 3).Value = 6
$ws.Cells.Item(38, 6).Value = 118
$ws.Cells.Item(38, 8).Value = 'kitchens'
$ws.Cells.Item(38, 9).Value = 'distractor'
$ws.Cells.Item(38, 11).Value = 'f'
$ws.Cells.Item(38, 12).Value = 'stimuli/img_uegbb.png'
$ws.Cells.Item(38, 13).Value = 78.80952380952381
$ws.Cells.Item(38, 14).Value = 61.52380952380953
$ws.Cells.Item(38, 15).Value = 70.16666666666667
$ws.Cells.Item(38, 16).Value = 42
$ws.Cells.Item(38, 17).Value = 8
$ws.Cells.Item(38, 18).Value = 8
$ws.Cells.Item(38, 19).Value = 8
$ws.Cells.Item(38, 20).Value = 8
$ws.Cells.Item(38, 21).Value = 8
$ws.Cells.Item(38, 22).Value = 8

# Row 39
$ws.Cells.Item(39, 3).Value = 6
$ws.Cells.Item(39, 6).Value = 119
$ws.Cells.Item(39, 8).Value = 'living_rooms'
$ws.Cells.Item(39, 9).Value = 'target'
$ws.Cells.Item(39, 11).Value = 'j'
$ws.Cells.Item(39, 12).Value = 'stimuli/img_3sw8t.png'
$ws.Cells.Item(39, 13).Value = 67.4888888888889
$ws.Cells.Item(39, 14).Value = 48.51111111111111
$ws.Cells.Item(39, 15).Value = 58
$ws.Cells.Item(39, 16).Value = 45
$ws.Cells.Item(39, 17).Value = 5
$ws.Cells.Item(39, 18).Value = 5
$ws.Cells.Item(39, 19).Value = 5
$ws.Cells.Item(39, 20).Value = 5
$ws.Cells.Item(39, 21).Value = 5
$ws.Cells.Item(39, 22).Value = 5

# Row 40
$ws.Cells.Item(40, 3).Value = 6
$ws.Cells.Item(40, 6).Value = 120
$ws.Cells.Item(40, 12).Value = 'stimuli/img_ub9nn.png'
$ws.Cells.Item(40, 13).Value = 78.77142857142857
$ws.Cells.Item(40, 14).Value = 60.37142857142857
$ws.Cells.Item(40, 15).Value = 69.57142857142857
$ws.Cells.Item(40, 16).Value = 35
$ws.Cells.Item(40, 17).Value = 8
$ws.Cells.Item(40, 18).Value = 8
$ws.Cells.Item(40, 19).Value = 8
$ws.Cells.Item(40, 20).Value = 8
$ws.Cells.Item(40, 21).Value = 8
$ws.Cells.Item(40, 22).Value = 8

# Row 41
$ws.Cells.Item(41, 3).Value = 6
$ws.Cells.Item(41, 6).Value = 121
$ws.Cells.Item(41, 8).Value = 'living_rooms'
$ws.Cells.Item(41, 9).Value = 'target'
$ws.Cells.Item(41, 11).Value = 'j'
$ws.Cells.Item(41, 12).Value = 'stimuli/img_pbsj1.png'
$ws.Cells.Item(41, 13).Value = 73.88636363636364
$ws.Cells.Item(41, 14).Value = 51.52272727272727
$ws.Cells.Item(41, 15).Value = 62.70454545454545
$ws.Cells.Item(41, 16).Value = 44
$ws.Cells.Item(41, 17).Value = 6
$ws.Cells.Item(41, 18).Value = 6
$ws.Cells.Item(41, 19).Value = 6
$ws.Cells.Item(41, 20).Value = 6
$ws.Cells.Item(41, 21).Value = 6
$ws.Cells.Item(41, 22).Value = 6
